$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -22.37940000000002
$ws.Range("B4").Value = 5.4847
$ws.Range("A6").Value = -22.90320000000001
$ws.Range("A7").Value = -19.98729999999998
$ws.Range("C7").Value = -12.49080000000001
$ws.Range("C8").Value = -11.98049999999999
$ws.Range("B9").Value = 6.188399999999997
$ws.Range("C10").Value = -13.7572
$ws.Range("B12").Value = 4.755599999999998
$ws.Range("C13").Value = -13.74599999999999
$ws.Range("A16").Value = -21.88209999999999
$ws.Range("C16").Value = -12.75890000000001
$ws.Range("B17").Value = 5.770499999999998
$ws.Range("B18").Value = 6.421599999999995
$ws.Range("B19").Value = 9.186899999999998
$ws.Range("A20").Value = -22.13040000000001
$ws.Range("B20").Value = 4.965899999999997
$ws.Range("B26").Value = 4.274200000000003
$ws.Range("A28").Value = -21.9289
$ws.Range("A29").Value = -21.32889999999998
$ws.Range("C30").Value = -11.9966
$ws.Range("B31").Value = 4.055599999999997
$ws.Range("A32").Value = -21.22990000000001
$ws.Range("B39").Value = 9.348300000000005
$ws.Range("A40").Value = -21.78629999999998
$ws.Range("B40").Value = 5.846300000000006
$ws.Range("C40").Value = -12.5063
$ws.Range("B41").Value = 9.318399999999992
$ws.Range("B42").Value = 9.394099999999993
$ws.Range("B43").Value = 6.296900000000004
$ws.Range("C44").Value = -13.42929999999999
$ws.Range("A46").Value = -21.7591
$ws.Range("B47").Value = 5.567500000000003
$ws.Range("B48").Value = 5.531000000000006
$ws.Range("A51").Value = -21.785
$ws.Range("A52").Value = -22.17149999999999
$ws.Range("A57").Value = -21.78430000000002
$ws.Range("A59").Value = -22.09580000000001
$ws.Range("A62").Value = -22.11800000000002
$ws.Range("B63").Value = 4.758499999999998
$ws.Range("B64").Value = 5.3964
$ws.Range("A66").Value = -21.42750000000001
$ws.Range("A73").Value = -20.08179999999998
$ws.Range("A74").Value = -21.56769999999997
$ws.Range("B76").Value = 5.348600000000001
$ws.Range("B81").Value = 5.094500000000003
$ws.Range("B89").Value = 5.623699999999996
$ws.Range("C89").Value = -13.78009999999999
$ws.Range("C91").Value = -12.64100000000001
$ws.Range("A92").Value = -21.63860000000001
$ws.Range("B94").Value = 4.761699999999991
$ws.Range("A100").Value = -22.12020000000001
